$d = $word.ActiveDocument

# Work on the last section's primary header (this also creates footnotes.xml /
# endnotes.xml / header1.xml and wires up the header/footer paragraph styles,
# just like Word does the first time a header is inserted in a document).
$section = $d.Sections.Last
$header = $section.Headers.Item(1)  # wdHeaderFooterPrimary = 1

$header.Range.Text = "Julian Motta Jennifer Castro Carlos Tafurt Dennis Masso"
$header.Range.Paragraphs.Item(1).Style = $d.Styles.Item("Encabezado")
